$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be parsed as numbers
$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Set new values
$ws.Range("D2").Value = '29.584.57'
$ws.Range("E2").Value = '  -2.55%  '
$ws.Range("D3").Value = '2.001.20'
$ws.Range("E3").Value = '  -4.83%  '
$ws.Range("D4").Value = '1.014'
$ws.Range("E4").Value = '  +0.47%  '
$ws.Range("D5").Value = '330.48'
$ws.Range("E5").Value = '  -4.02%  '
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").Value = '0.5005'
$ws.Range("E7").Value = '  -4.16%  '
$ws.Range("D8").Value = '0.4239'
$ws.Range("E8").Value = '  -4.17%  '
$ws.Range("D9").Value = '54.19'
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").Value = '0.08961'
$ws.Range("E10").Value = '  -4.46%  '
$ws.Range("D11").Value = '1.121'
$ws.Range("E11").Value = '  -4.19%  '
$ws.Range("D12").Value = '23.38'
$ws.Range("E12").Value = '  -5.72%  '
$ws.Range("D13").Value = '2.007.10'
$ws.Range("E13").Value = '  -6.40%  '
$ws.Range("D14").Value = '8.070'
$ws.Range("E14").Value = '  -6.67%  '
$ws.Range("D15").Value = '6.497'
$ws.Range("E15").Value = '  -6.00%  '
$ws.Range("D16").Value = '1.014'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '94.46'
$ws.Range("E17").Value = '  -7.04%  '
$ws.Range("D18").Value = '0.00001114'
$ws.Range("E18").Value = '  -3.90%  '
$ws.Range("D19").Value = '0.06676'
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("D20").Value = '19.77'
$ws.Range("E20").Value = '  -6.48%  '
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("D22").Value = '5.969'
$ws.Range("E22").Value = '  -6.27%  '
$ws.Range("D23").Value = '29.593.57'
$ws.Range("E23").Value = '  -2.72%  '
$ws.Range("D24").Value = '11.98'
$ws.Range("E24").Value = '  -4.44%  '
$ws.Range("E25").Value = '  -0.90%  '
$ws.Range("D26").Value = '159.63'
$ws.Range("E26").Value = '  -1.74%  '
$ws.Range("D27").Value = '20.71'
$ws.Range("E27").Value = '  -5.35%  '
$ws.Range("D28").Value = '6.340'
$ws.Range("E28").Value = '  -5.80%  '
$ws.Range("D29").Value = '2.306'
$ws.Range("E29").Value = '  -8.18%  '
$ws.Range("D30").Value = '128.54'
$ws.Range("E30").Value = '  -3.67%  '
$ws.Range("D31").Value = '1.057'
$ws.Range("E31").Value = '  -6.74%  '
$ws.Range("D32").Value = '0.09954'
$ws.Range("E32").Value = '  -5.31%  '
$ws.Range("D33").Value = '1.566'
$ws.Range("E33").Value = '  -5.82%  '
$ws.Range("D34").Value = '5.851'
$ws.Range("E34").Value = '  -5.96%  '
$ws.Range("D35").Value = '3.784'
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("D36").Value = '9.444'
$ws.Range("E36").Value = '  -8.36%  '
$ws.Range("E37").Value = '  -6.13%  '
$ws.Range("D38").Value = '1.309'
$ws.Range("E38").Value = '  -2.62%  '
$ws.Range("D39").Value = '0.06359'
$ws.Range("E39").Value = '  -6.04%  '
$ws.Range("D40").Value = '0.6576'
$ws.Range("E40").Value = '  -6.24%  '
$ws.Range("D41").Value = '11.69'
$ws.Range("E41").Value = '  -6.58%  '
$ws.Range("D42").Value = '0.2056'
$ws.Range("E42").Value = '  -7.44%  '
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("D44").Value = '0.6339'
$ws.Range("E44").Value = '  -7.09%  '
$ws.Range("D45").Value = '13.52'
$ws.Range("E45").Value = '  -6.66%  '
$ws.Range("E46").Value = '  -5.97%  '
$ws.Range("D47").Value = '1.314'
$ws.Range("E47").Value = '  -5.96%  '
$ws.Range("D48").Value = '3.524'
$ws.Range("E48").Value = '  -3.18%  '
$ws.Range("D49").Value = '0.00000000340'
$ws.Range("E49").Value = '  -2.41%  '
$ws.Range("D50").Value = '0.06992'
$ws.Range("E50").Value = '  -3.56%  '
$ws.Range("D51").Value = '1.124'
$ws.Range("E51").Value = '  -7.29%  '

# Restore default style (style index 0) for cells we temporarily reformatted
foreach ($addr in $textCells) {
  $ws.Range($addr).Style = "Normal"
}
